$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect numeric-looking text in column D from Excel auto-converting
# the literal strings to actual numbers: force Text format first,
# then strip the format residue after all values are written so the
# cells end up as plain (unstyled) text cells, matching the source.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.342.70"
$ws.Range("E2").Value = "  -0.74%  "

$ws.Range("D3").Value = "1.711.31"
$ws.Range("E3").Value = "  -0.66%  "

$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "224.57"
$ws.Range("E5").Value = "  -0.40%  "

$ws.Range("D6").Value = "0.5304"

$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "0.06696"
$ws.Range("E8").Value = "  +1.46%  "

$ws.Range("D9").Value = "0.2669"

$ws.Range("E10").Value = "  -3.70%  "

$ws.Range("D11").Value = "0.07681"
$ws.Range("E11").Value = "  -0.46%  "

$ws.Range("D12").Value = "4.517"
$ws.Range("E12").Value = "  -2.08%  "

$ws.Range("D13").Value = "1.946.76"
$ws.Range("E13").Value = "  -0.69%  "

$ws.Range("D14").Value = "1.712.31"
$ws.Range("E14").Value = "  -0.72%  "

$ws.Range("D15").Value = "0.5843"
$ws.Range("E15").Value = "  -0.05%  "

$ws.Range("D16").Value = "0.0₅8230"
$ws.Range("E16").Value = "  -0.92%  "

$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("D18").Value = "27.363.09"
$ws.Range("E18").Value = "  -0.71%  "

$ws.Range("D19").Value = "223.38"
$ws.Range("E19").Value = "  +1.26%  "

$ws.Range("E20").Value = "  +0.03%  "

$ws.Range("D21").Value = "4.640"
$ws.Range("E21").Value = "  -1.79%  "

$ws.Range("D22").Value = "10.41"
$ws.Range("E22").Value = "  -2.15%  "

$ws.Range("D23").Value = "6.007"
$ws.Range("E23").Value = "  -1.35%  "

$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("D25").Value = "144.40"
$ws.Range("E25").Value = "  -2.78%  "

$ws.Range("D26").Value = "1.693"
$ws.Range("E26").Value = "  -2.66%  "

$ws.Range("D27").Value = "0.1214"
$ws.Range("E27").Value = "  -1.68%  "

$ws.Range("D28").Value = "7.251"
$ws.Range("E28").Value = "  -2.15%  "

$ws.Range("D29").Value = "16.26"
$ws.Range("E29").Value = "  -2.32%  "

$ws.Range("D30").Value = "0.05373"
$ws.Range("E30").Value = "  -4.03%  "

$ws.Range("D31").Value = "1.293"
$ws.Range("E31").Value = "  -0.73%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "3.442"
$ws.Range("E32").Value = "  -0.31%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "3.456"
$ws.Range("E33").Value = "  -2.92%  "

$ws.Range("D34").Value = "1.643"
$ws.Range("E34").Value = "  -1.03%  "

$ws.Range("D35").Value = "2.865"
$ws.Range("E35").Value = "  +1.03%  "

$ws.Range("D36").Value = "0.9523"
$ws.Range("E36").Value = "  -0.87%  "

$ws.Range("D37").Value = "2.390"
$ws.Range("E37").Value = "  -1.52%  "

$ws.Range("D38").Value = "0.5874"
$ws.Range("E38").Value = "  -1.33%  "

$ws.Range("E39").Value = "  -0.78%  "

$ws.Range("D40").Value = "1.092.87"
$ws.Range("E40").Value = "  +3.81%  "

$ws.Range("D41").Value = "5.812"
$ws.Range("E41").Value = "  -1.80%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "0.8442"
$ws.Range("E42").Value = "  -1.10%  "

$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "1.004"
$ws.Range("E43").Value = "  +0.07%  "

$ws.Range("D44").Value = "100.98"
$ws.Range("E44").Value = "  -0.33%  "

$ws.Range("D45").Value = "1.854.42"
$ws.Range("E45").Value = "  -0.70%  "

$ws.Range("D46").Value = "0.0₈115"
$ws.Range("E46").Value = "  -0.22%  "

$ws.Range("D47").Value = "57.95"
$ws.Range("E47").Value = "  -1.88%  "

$ws.Range("D48").Value = "0.4535"
$ws.Range("E48").Value = "  +2.25%  "

$ws.Range("D49").Value = "1.011"
$ws.Range("E49").Value = "  +1.00%  "

$ws.Range("D50").Value = "8.128"
$ws.Range("E50").Value = "  -0.78%  "

$ws.Range("D51").Value = "0.05235"
$ws.Range("E51").Value = "  -0.35%  "

# Remove the temporary Text-format styling so D-column cells revert
# to the default (unstyled) appearance, same as the rest of the sheet.
$ws.Range("D2:D51").ClearFormats()
